$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308, shifting existing rows 308:352 down to 309:353
$ws.Rows("308:308").Insert()

# Populate the newly inserted row 308 with its data
$ws.Cells.Item(308, 1).Value = 10
$ws.Cells.Item(308, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(308, 3).Value = "La Araucanía"
$ws.Cells.Item(308, 4).Value = (Get-Date -Year 2023 -Month 2 -Day 27 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(308, 5).Value = 9
$ws.Cells.Item(308, 6).Value = 100112052
$ws.Cells.Item(308, 7).Value = "Albahaca"
$ws.Cells.Item(308, 8).Value = "Sin especificar"
$ws.Cells.Item(308, 9).Value = "Primera"
$ws.Cells.Item(308, 10).Value = 40
$ws.Cells.Item(308, 11).Value = 6000
$ws.Cells.Item(308, 12).Value = 6000
$ws.Cells.Item(308, 13).Value = 6000
$ws.Cells.Item(308, 14).Value = "$/paquete"
$ws.Cells.Item(308, 15).Value = "Región del Maule"
$ws.Cells.Item(308, 16).Value = 6000
$ws.Cells.Item(308, 17).Value = 1
$ws.Cells.Item(308, 18).Value = "Hortaliza"
